$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function FindRowByAccount($account) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($i = 1; $i -le $lastRow; $i++) {
        $acct = $ws.Cells.Item($i, 1).Value2
        if ($acct -eq $account) {
            return $i
        }
    }
    return -1
}

# --- Step 1: remove the four rows that disappeared from the top of the
# ledger (KARINA, EDUARDO, PEDRO, LEDA) ---
$karinaRow = FindRowByAccount("005141215")
$ws.Rows("$karinaRow`:$($karinaRow + 3)").Delete()

# --- Step 2: insert the new AYRTON row in their place with the new
# balance, right after CARLOS (row 2) ---
$carlosRow = FindRowByAccount("004211922")
$newRow = $carlosRow + 1
$ws.Rows($newRow).Insert()
$acctCell = $ws.Cells.Item($newRow, 1)
$acctCell.NumberFormat = "@"
$acctCell.Value = "001000882"
$ws.Cells.Item($newRow, 2).Value = "AYRTON"
$ws.Cells.Item($newRow, 3).Value = 9918.09

# --- Step 3: update CHRISTIAN's balance ---
$christianRow = FindRowByAccount("004420763")
$ws.Cells.Item($christianRow, 3).Value = 8000

# --- Step 4: insert a new HEITOR row right before GUSTAVO (004313254) ---
$gustavoRow = FindRowByAccount("004313254")
$ws.Rows($gustavoRow).Insert()
$acctCell2 = $ws.Cells.Item($gustavoRow, 1)
$acctCell2.NumberFormat = "@"
$acctCell2.Value = "003435941"
$ws.Cells.Item($gustavoRow, 2).Value = "HEITOR"
$ws.Cells.Item($gustavoRow, 3).Value = 807.5

# --- Step 5: remove the old AYRTON row (93.41) further down the sheet ---
$oldAyrtonRow = -1
$lastRow = $ws.UsedRange.Rows.Count
for ($i = 1; $i -le $lastRow; $i++) {
    $acct = $ws.Cells.Item($i, 1).Value2
    $name = $ws.Cells.Item($i, 2).Value2
    if (($acct -eq "001000882") -and ($name -eq "AYRTON") -and ($i -ne $newRow)) {
        $oldAyrtonRow = $i
        break
    }
}
$ws.Rows($oldAyrtonRow).Delete()
